$wb = $excel.ActiveWorkbook

# --- "runs" sheet: Start/Stop Run ID + Sim Time ---
$runs = $wb.Worksheets.Item("runs")
$runs.Range("B2").Value = 15
$runs.Range("B3").Value = 500

# --- "params" sheet: per-run configuration table ---
$params = $wb.Worksheets.Item("params")

$params.Range("B2").Value = 10
$params.Range("C2").Value = 0
$params.Range("D2").Value = 1
$params.Range("J2").Value = 250
$params.Range("K2").Value = "config1"

$params.Range("B3").Value = 20
$params.Range("C3").Value = 0
$params.Range("D3").Value = 0.9
$params.Range("J3").Value = 250
$params.Range("K3").Value = "config2"

$params.Range("B4").Value = 30
$params.Range("C4").Value = 0
$params.Range("D4").Value = 0.8
$params.Range("J4").Value = 250
$params.Range("K4").Value = "config3"

$params.Range("B5").Value = 40
$params.Range("C5").Value = 0
$params.Range("D5").Value = 0.7
$params.Range("J5").Value = 250
$params.Range("K5").Value = "config4"

$params.Range("B6").Value = 50
$params.Range("C6").Value = 0
$params.Range("D6").Value = 0
$params.Range("J6").Value = 250
$params.Range("K6").Value = "config5"

$params.Range("B7").Value = 0
$params.Range("C7").Value = 10
$params.Range("D7").Value = 1
$params.Range("J7").Value = 250
$params.Range("K7").Value = "config1"

$params.Range("B8").Value = 0
$params.Range("C8").Value = 20
$params.Range("D8").Value = 0.9
$params.Range("J8").Value = 250
$params.Range("K8").Value = "config2"

$params.Range("B9").Value = 0
$params.Range("C9").Value = 30
$params.Range("D9").Value = 0.8
$params.Range("J9").Value = 250
$params.Range("K9").Value = "config3"

$params.Range("B10").Value = 0
$params.Range("C10").Value = 40
$params.Range("D10").Value = 0.7
$params.Range("J10").Value = 250
$params.Range("K10").Value = "config4"

$params.Range("B11").Value = 0
$params.Range("C11").Value = 50
$params.Range("D11").Value = 0
$params.Range("J11").Value = 250
$params.Range("K11").Value = "config5"

$params.Range("B12").Value = 10
$params.Range("C12").Value = 0
$params.Range("D12").Value = 1
$params.Range("J12").Value = 250
$params.Range("K12").Value = "config1"

$params.Range("B13").Value = 10
$params.Range("C13").Value = 10
$params.Range("D13").Value = 0.9
$params.Range("J13").Value = 250
$params.Range("K13").Value = "config2"

$params.Range("B14").Value = 15
$params.Range("C14").Value = 15
$params.Range("D14").Value = 0.8
$params.Range("J14").Value = 250
$params.Range("K14").Value = "config3"

$params.Range("B15").Value = 20
$params.Range("C15").Value = 20
$params.Range("D15").Value = 0.7
$params.Range("J15").Value = 250
$params.Range("K15").Value = "config4"

$params.Range("B16").Value = 25
$params.Range("C16").Value = 25
$params.Range("D16").Value = 0
$params.Range("J16").Value = 250
$params.Range("K16").Value = "config5"

# Match the author's final selection state on the "params" sheet.
$null = $params.Range("J2:J16").Select()

Write-Host "edits applied"
